$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9242176413536072
$ws.Range("B1").Value = 2.028493642807007
$ws.Range("C1").Value = 8.804119110107422
$ws.Range("D1").Value = 1.833012342453003
$ws.Range("E1").Value = 1.426158785820007
